# Reorders the "Requisitos" bullet list in the LOQ4230 document, drops the
# "LOB1045 - Leitura e Producao de Textos Academicos" requisite, and fixes
# the "Algebra" -> "Álgebra" accent typo.
#
# The paragraph is a single w:p holding one w:r run per requisite
# ("<w:t>TEXT</w:t><w:br/>"). A handful of plain Range.Text edits would
# normally do the trick, but this host coalesces any two neighbouring runs
# that both get touched (same rPr) into a single w:r, which would collapse
# the whole list into one run and not match the expected per-item run
# layout. InsertXML lets us hand the host fully-formed markup instead, so
# the new list keeps one run per requisite, exactly like before the edit.

$d = $word.ActiveDocument

# New order for the list, top to bottom, with the accent fixed and the
# "LOB1045" item removed.
$newItems = @(
    "LOQ4251 -  Fundamentos de Química  (Requisito)",
    "LOB1006 -  Cálculo IV  (Requisito)",
    "LOB1053 -  Física III  (Requisito)",
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)",
    "LOB1003 -  Cálculo I  (Requisito)",
    "LOB1012 -  Estatística  (Requisito)",
    "LOB1024 -  Mecânica  (Requisito)",
    "LOB1036 -  Geometria Analítica  (Requisito)",
    "LOB1037 -  Álgebra Linear  (Requisito)",
    "LOB1038 -  Física Experimental I  (Requisito)",
    "LOB1039 -  Física Experimental III  (Requisito)",
    "LOB1041 -  Física Experimental II  (Requisito)",
    "LOB1052 -  Cálculo III  (Requisito)",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
    "LOB1004 -  Cálculo II  (Requisito)",
    "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)",
    "LOB1018 -  Física I  (Requisito)",
    "LOB1019 -  Física II  (Requisito)"
)

# Locate the list paragraph via the (still untouched) first requisite's
# text, rather than hard-coding a paragraph index.
$locator = $d.Content
[void]$locator.Find.Execute("LOB1003 -  Cálculo I  (Requisito)")
$hitStart = $locator.Start

$paraStart = -1
$paraEnd = -1
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -le $hitStart -and $para.Range.End -gt $hitStart) {
        $paraStart = $para.Range.Start
        $paraEnd = $para.Range.End
        break
    }
}

# Build the replacement paragraph XML: same ListBullet style, one w:r per
# requisite (text + line break), in the new order.
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$sb = New-Object System.Text.StringBuilder
[void]$sb.Append("<w:p $ns><w:pPr><w:pStyle w:val=""ListBullet""/></w:pPr>")
foreach ($item in $newItems) {
    $escaped = $item.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    [void]$sb.Append("<w:r><w:t>$escaped</w:t><w:br/></w:r>")
}
[void]$sb.Append("</w:p>")
$xml = $sb.ToString()

$full = $d.Range($paraStart, $paraEnd - 1)
[void]$full.InsertXML($xml)
